$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46
$ws.Range("A44:G44").Copy()
$ws.Range("A46:G46").PasteSpecial(-4122)
$ws.Cells.Item(46,1).Value = "Year 2"
$ws.Cells.Item(46,2).Value = "C1"
$ws.Cells.Item(46,3).Value = "anatomy"
$ws.Cells.Item(46,4).Formula = '="3"'
$ws.Cells.Item(46,5).Formula = '="13/10/2025"'
$ws.Cells.Item(46,6).Formula = '="12:00:00"'
$ws.Cells.Item(46,7).Value = 75
$ws.Range("D46:F46").Copy()
$ws.Range("D46:F46").PasteSpecial(-4163)

# Row 47
$ws.Range("A45:G45").Copy()
$ws.Range("A47:G47").PasteSpecial(-4122)
$ws.Cells.Item(47,1).Value = "Year 2"
$ws.Cells.Item(47,2).Value = "C1"
$ws.Cells.Item(47,3).Value = "biochemistry cbl"
$ws.Cells.Item(47,4).Formula = '="1"'
$ws.Cells.Item(47,5).Formula = '="30/09/2025"'
$ws.Cells.Item(47,6).Formula = '="12:00:00"'
$ws.Cells.Item(47,7).Value = 120
$ws.Range("D47:F47").Copy()
$ws.Range("D47:F47").PasteSpecial(-4163)

# Row 48
$ws.Range("A44:G44").Copy()
$ws.Range("A48:G48").PasteSpecial(-4122)
$ws.Cells.Item(48,1).Value = "Year 2"
$ws.Cells.Item(48,2).Value = "C1"
$ws.Cells.Item(48,3).Value = "histology"
$ws.Cells.Item(48,4).Formula = '="1"'
$ws.Cells.Item(48,5).Formula = '="30/09/2025"'
$ws.Cells.Item(48,6).Formula = '="14:00:00"'
$ws.Cells.Item(48,7).Value = 120
$ws.Range("D48:F48").Copy()
$ws.Range("D48:F48").PasteSpecial(-4163)

# Row 49
$ws.Range("A45:G45").Copy()
$ws.Range("A49:G49").PasteSpecial(-4122)
$ws.Cells.Item(49,1).Value = "Year 2"
$ws.Cells.Item(49,2).Value = "C1"
$ws.Cells.Item(49,3).Value = "histology"
$ws.Cells.Item(49,4).Formula = '="2"'
$ws.Cells.Item(49,5).Formula = '="06/10/2025"'
$ws.Cells.Item(49,6).Formula = '="12:00:00"'
$ws.Cells.Item(49,7).Value = 120
$ws.Range("D49:F49").Copy()
$ws.Range("D49:F49").PasteSpecial(-4163)

# Row 50
$ws.Range("A44:G44").Copy()
$ws.Range("A50:G50").PasteSpecial(-4122)
$ws.Cells.Item(50,1).Value = "Year 2"
$ws.Cells.Item(50,2).Value = "C1"
$ws.Cells.Item(50,3).Value = "histology"
$ws.Cells.Item(50,4).Formula = '="3"'
$ws.Cells.Item(50,5).Formula = '="13/10/2025"'
$ws.Cells.Item(50,6).Formula = '="10:00:00"'
$ws.Cells.Item(50,7).Value = 120
$ws.Range("D50:F50").Copy()
$ws.Range("D50:F50").PasteSpecial(-4163)

# Row 51
$ws.Range("A45:G45").Copy()
$ws.Range("A51:G51").PasteSpecial(-4122)
$ws.Cells.Item(51,1).Value = "Year 2"
$ws.Cells.Item(51,2).Value = "C1"
$ws.Cells.Item(51,3).Value = "microbiology"
$ws.Cells.Item(51,4).Formula = '="3"'
$ws.Cells.Item(51,5).Formula = '="14/10/2025"'
$ws.Cells.Item(51,6).Formula = '="14:00:00"'
$ws.Cells.Item(51,7).Value = 120
$ws.Range("D51:F51").Copy()
$ws.Range("D51:F51").PasteSpecial(-4163)

# Row 52
$ws.Range("A44:G44").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$ws.Cells.Item(52,1).Value = "Year 2"
$ws.Cells.Item(52,2).Value = "C1"
$ws.Cells.Item(52,3).Value = "parasitology"
$ws.Cells.Item(52,4).Formula = '="1"'
$ws.Cells.Item(52,5).Formula = '="29/09/2025"'
$ws.Cells.Item(52,6).Formula = '="08:00:00"'
$ws.Cells.Item(52,7).Value = 120
$ws.Range("D52:F52").Copy()
$ws.Range("D52:F52").PasteSpecial(-4163)

# Row 53
$ws.Range("A45:G45").Copy()
$ws.Range("A53:G53").PasteSpecial(-4122)
$ws.Cells.Item(53,1).Value = "Year 2"
$ws.Cells.Item(53,2).Value = "C1"
$ws.Cells.Item(53,3).Value = "parasitology"
$ws.Cells.Item(53,4).Formula = '="2"'
$ws.Cells.Item(53,5).Formula = '="06/10/2025"'
$ws.Cells.Item(53,6).Formula = '="14:00:00"'
$ws.Cells.Item(53,7).Value = 120
$ws.Range("D53:F53").Copy()
$ws.Range("D53:F53").PasteSpecial(-4163)

# Row 54
$ws.Range("A44:G44").Copy()
$ws.Range("A54:G54").PasteSpecial(-4122)
$ws.Cells.Item(54,1).Value = "Year 2"
$ws.Cells.Item(54,2).Value = "C1"
$ws.Cells.Item(54,3).Value = "parasitology"
$ws.Cells.Item(54,4).Formula = '="2"'
$ws.Cells.Item(54,5).Formula = '="12/10/2025"'
$ws.Cells.Item(54,6).Formula = '="12:00:00"'
$ws.Cells.Item(54,7).Value = 120
$ws.Range("D54:F54").Copy()
$ws.Range("D54:F54").PasteSpecial(-4163)

# Row 55
$ws.Range("A45:G45").Copy()
$ws.Range("A55:G55").PasteSpecial(-4122)
$ws.Cells.Item(55,1).Value = "Year 2"
$ws.Cells.Item(55,2).Value = "C1"
$ws.Cells.Item(55,3).Value = "parasitology"
$ws.Cells.Item(55,4).Formula = '="3"'
$ws.Cells.Item(55,5).Formula = '="13/10/2025"'
$ws.Cells.Item(55,6).Formula = '="14:00:00"'
$ws.Cells.Item(55,7).Value = 120
$ws.Range("D55:F55").Copy()
$ws.Range("D55:F55").PasteSpecial(-4163)

# Row 56
$ws.Range("A44:G44").Copy()
$ws.Range("A56:G56").PasteSpecial(-4122)
$ws.Cells.Item(56,1).Value = "Year 2"
$ws.Cells.Item(56,2).Value = "C1"
$ws.Cells.Item(56,3).Value = "pathology lab"
$ws.Cells.Item(56,4).Formula = '="2"'
$ws.Cells.Item(56,5).Formula = '="13/10/2025"'
$ws.Cells.Item(56,6).Formula = '="08:00:00"'
$ws.Cells.Item(56,7).Value = 120
$ws.Range("D56:F56").Copy()
$ws.Range("D56:F56").PasteSpecial(-4163)

# Row 57
$ws.Range("A45:G45").Copy()
$ws.Range("A57:G57").PasteSpecial(-4122)
$ws.Cells.Item(57,1).Value = "Year 2"
$ws.Cells.Item(57,2).Value = "C1"
$ws.Cells.Item(57,3).Value = "pharmacology"
$ws.Cells.Item(57,4).Formula = '="1"'
$ws.Cells.Item(57,5).Formula = '="01/10/2025"'
$ws.Cells.Item(57,6).Formula = '="08:00:00"'
$ws.Cells.Item(57,7).Value = 120
$ws.Range("D57:F57").Copy()
$ws.Range("D57:F57").PasteSpecial(-4163)

# Row 58
$ws.Range("A44:G44").Copy()
$ws.Range("A58:G58").PasteSpecial(-4122)
$ws.Cells.Item(58,1).Value = "Year 2"
$ws.Cells.Item(58,2).Value = "C1"
$ws.Cells.Item(58,3).Value = "pharmacology"
$ws.Cells.Item(58,4).Formula = '="2"'
$ws.Cells.Item(58,5).Formula = '="05/10/2025"'
$ws.Cells.Item(58,6).Formula = '="08:00:00"'
$ws.Cells.Item(58,7).Value = 120
$ws.Range("D58:F58").Copy()
$ws.Range("D58:F58").PasteSpecial(-4163)

# Row 59
$ws.Range("A45:G45").Copy()
$ws.Range("A59:G59").PasteSpecial(-4122)
$ws.Cells.Item(59,1).Value = "Year 2"
$ws.Cells.Item(59,2).Value = "C1"
$ws.Cells.Item(59,3).Value = "pharmacology"
$ws.Cells.Item(59,4).Formula = '="3"'
$ws.Cells.Item(59,5).Formula = '="14/10/2025"'
$ws.Cells.Item(59,6).Formula = '="12:00:00"'
$ws.Cells.Item(59,7).Value = 120
$ws.Range("D59:F59").Copy()
$ws.Range("D59:F59").PasteSpecial(-4163)

# Row 60
$ws.Range("A44:G44").Copy()
$ws.Range("A60:G60").PasteSpecial(-4122)
$ws.Cells.Item(60,1).Value = "Year 2"
$ws.Cells.Item(60,2).Value = "C1"
$ws.Cells.Item(60,3).Value = "physiology"
$ws.Cells.Item(60,4).Formula = '="2"'
$ws.Cells.Item(60,5).Formula = '="07/10/2025"'
$ws.Cells.Item(60,6).Formula = '="08:00:00"'
$ws.Cells.Item(60,7).Value = 120
$ws.Range("D60:F60").Copy()
$ws.Range("D60:F60").PasteSpecial(-4163)

$excel.CutCopyMode = 0
